$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("Анджигаев Владислав"): raise scores in columns G-J from 2 to 5.
# L5 = SUM(C5:J5) recalculates automatically from 28 to 40.
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 5

# Move the active selection to I5 (also updates the frozen-pane view state).
$ws.Range("I5").Select()
